# Updates the cryptos.xlsx price/volume snapshot (GitHub Actions refresh).
# Cells are written as text (NumberFormat "@") and the style is reset back
# to "Normal" afterwards so numeric-looking strings (e.g. "0.111", "1.00")
# are preserved verbatim instead of being coerced into numbers by Excel,
# and so no lingering custom number format is left applied to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws 'D2' '62.914.04'
Set-TextCell $ws 'E2' '  +0.34%  '
Set-TextCell $ws 'D3' '2.463.33'
Set-TextCell $ws 'E3' '  +0.78%  '
Set-TextCell $ws 'D5' '575.34'
Set-TextCell $ws 'E5' '  -0.07%  '
Set-TextCell $ws 'D6' '146.62'
Set-TextCell $ws 'E7' '  -0.02%  '
Set-TextCell $ws 'E8' '  -0.12%  '
Set-TextCell $ws 'D9' '2.463.55'
Set-TextCell $ws 'E9' '  +0.79%  '
Set-TextCell $ws 'D10' '0.111'
Set-TextCell $ws 'E10' '  +1.28%  '
Set-TextCell $ws 'D11' '0.162'
Set-TextCell $ws 'E11' '  +1.26%  '
Set-TextCell $ws 'E12' '  +0.73%  '
Set-TextCell $ws 'E13' '  +1.06%  '
Set-TextCell $ws 'D14' '29.04'
Set-TextCell $ws 'E14' '  +2.24%  '
Set-TextCell $ws 'E15' '  +0.24%  '
Set-TextCell $ws 'D16' '2.910.13'
Set-TextCell $ws 'E16' '  +0.75%  '
Set-TextCell $ws 'D17' '62.798.20'
Set-TextCell $ws 'E17' '  +0.35%  '
Set-TextCell $ws 'D18' '2.463.80'
Set-TextCell $ws 'E18' '  +0.78%  '
Set-TextCell $ws 'E19' '  +2.65%  '
Set-TextCell $ws 'D20' '11.03'
Set-TextCell $ws 'E20' '  +1.10%  '
Set-TextCell $ws 'D21' '327.36'
Set-TextCell $ws 'E21' '  +0.21%  '
Set-TextCell $ws 'D22' '2.25'
Set-TextCell $ws 'E22' '  +11.92%  '
Set-TextCell $ws 'E23' '  +0.00%  '
Set-TextCell $ws 'B24' 'Dai'
Set-TextCell $ws 'C24' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws 'D24' '1.00'
Set-TextCell $ws 'E24' '  +0.05%  '
Set-TextCell $ws 'B25' 'Aptos'
Set-TextCell $ws 'C25' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D25' '10.27'
Set-TextCell $ws 'E25' '  +21.21%  '
Set-TextCell $ws 'D26' '65.85'
Set-TextCell $ws 'E26' '  +0.48%  '
Set-TextCell $ws 'D27' '652.26'
Set-TextCell $ws 'D28' ('0.0{0}0984' -f [char]0x2083)
Set-TextCell $ws 'E28' '  +0.71%  '
Set-TextCell $ws 'D29' '2.583.88'
Set-TextCell $ws 'E30' '  -13.92%  '
Set-TextCell $ws 'E31' '  +2.52%  '
Set-TextCell $ws 'D32' '8.01'
Set-TextCell $ws 'E32' '  -1.90%  '
Set-TextCell $ws 'E33' '  -0.55%  '
Set-TextCell $ws 'E34' '  -3.72%  '
Set-TextCell $ws 'D35' '0.999'
Set-TextCell $ws 'E35' '  +0.00%  '
Set-TextCell $ws 'E36' '  +3.72%  '
Set-TextCell $ws 'E37' '  +0.66%  '
Set-TextCell $ws 'E38' '  -0.82%  '
Set-TextCell $ws 'B39' 'RenderToken'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextCell $ws 'D39' '5.42'
Set-TextCell $ws 'E39' '  -0.72%  '
Set-TextCell $ws 'B40' 'EthereumClassic'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws 'D40' '18.74'
Set-TextCell $ws 'E40' '  +0.81%  '
Set-TextCell $ws 'B41' 'Monero'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D41' '151.46'
Set-TextCell $ws 'E41' '  -1.13%  '
Set-TextCell $ws 'E42' '  +3.25%  '
Set-TextCell $ws 'D43' '1.73'
Set-TextCell $ws 'E43' '  -1.39%  '
Set-TextCell $ws 'D44' ('0.0{0}0319' -f [char]0x2086)
Set-TextCell $ws 'E44' '  -62.78%  '
Set-TextCell $ws 'E45' '  +0.00%  '
Set-TextCell $ws 'D46' '153.59'
Set-TextCell $ws 'E46' '  +6.56%  '
Set-TextCell $ws 'D47' '15.24'
Set-TextCell $ws 'E48' '  +0.05%  '
Set-TextCell $ws 'D49' '20.47'
Set-TextCell $ws 'E49' '  -0.47%  '
Set-TextCell $ws 'D50' '0.609'
Set-TextCell $ws 'E51' '  +0.06%  '
